# Apply cryptos list update (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.482.89"

$ws.Range("D3").Value = "1.637.42"
$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.88"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3785"
$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.77"
$ws.Range("E8").Value = "  -0.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3620"
$ws.Range("E9").Value = "  -0.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08162"
$ws.Range("E10").Value = "  +0.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.229"
$ws.Range("E11").Value = "  -1.92%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.47"
$ws.Range("E13").Value = "  -2.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.456"
$ws.Range("E14").Value = "  -3.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.352"
$ws.Range("E15").Value = "  +0.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001238"
$ws.Range("E16").Value = "  -1.85%  "

$ws.Range("D17").Value = "1.632.68"
$ws.Range("E17").Value = "  -0.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.46"
$ws.Range("E18").Value = "  +1.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06934"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.584"
$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.51"
$ws.Range("E21").Value = "  -3.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.49"
$ws.Range("E23").Value = "  -3.12%  "

$ws.Range("D24").Value = "23.470.87"
$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.498"
$ws.Range("E25").Value = "  +2.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.050"
$ws.Range("E26").Value = "  -6.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.14"
$ws.Range("E27").Value = "  -0.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.61"
$ws.Range("E28").Value = "  +0.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.249"
$ws.Range("E29").Value = "  -1.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.74"
$ws.Range("E30").Value = "  -2.67%  "

$ws.Range("D31").Value = "1.814.92"
$ws.Range("E31").Value = "  -0.66%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.088"
$ws.Range("E32").Value = "  +12.61%  "

$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.159"
$ws.Range("E33").Value = "  -7.31%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.585"
$ws.Range("E34").Value = "  -5.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.42"
$ws.Range("E35").Value = "  +3.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02758"
$ws.Range("E36").Value = "  -4.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2494"
$ws.Range("E37").Value = "  -2.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08772"
$ws.Range("E38").Value = "  -0.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07091"
$ws.Range("E39").Value = "  -2.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.983"
$ws.Range("E40").Value = "  -4.95%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.353"
$ws.Range("E41").Value = "  -2.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7049"
$ws.Range("E42").Value = "  -1.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.18"
$ws.Range("E43").Value = "  -3.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.55"
$ws.Range("E44").Value = "  -5.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6524"
$ws.Range("E45").Value = "  -0.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.274"
$ws.Range("E47").Value = "  -4.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.968"
$ws.Range("E48").Value = "  -0.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07968"
$ws.Range("E49").Value = "  -0.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.64"
$ws.Range("E50").Value = "  +0.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.191"
$ws.Range("E51").Value = "  -2.25%  "

